$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray x1 / Character1 values in rows 4 and 6 (columns H and I),
# leaving the cells blank but keeping their existing style/formatting.
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("I6").Value = $null

# Update the active selection to E24 (as last left by the author).
$ws.Range("E24").Select()
